$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 10912.5
$ws.Range("I43").Value = 9250
$ws.Range("J43").Value = 12575
$ws.Range("K43").Value = 9250
$ws.Range("L43").Value = 12575
$ws.Range("M43").Value = -9181
$ws.Range("N43").Value = -12713
$ws.Range("H68").Value = 25295
$ws.Range("J68").Value = 25295
$ws.Range("L68").Value = 25295
$ws.Range("N68").Value = -26793
$ws.Range("H71").Value = 25295
$ws.Range("J71").Value = 25295
$ws.Range("L71").Value = 75885
$ws.Range("N71").Value = -83373

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2066.6858
$ws.Range("I2").Value = 1914.96
$ws.Range("J2").Value = 2446
$ws.Range("K2").Value = 1914.96
$ws.Range("L2").Value = 2446
$ws.Range("M2").Value = -1801.96
$ws.Range("N2").Value = -2672
$ws.Range("H61").Value = 1467
$ws.Range("I61").Value = 1188.2727
$ws.Range("K61").Value = 1188.2727
$ws.Range("M61").Value = -976.2727
$ws.Range("H74").Value = 1233.325
$ws.Range("I74").Value = 1206.8286
$ws.Range("K74").Value = 1206.8286
$ws.Range("M74").Value = -332.8286000000001
$ws.Range("H77").Value = 1233.325
$ws.Range("I77").Value = 1206.8286
$ws.Range("K77").Value = 6034.143
$ws.Range("M77").Value = -1666.143
$ws.Range("H95").Value = 30208
$ws.Range("J95").Value = 30208
$ws.Range("L95").Value = 30208
$ws.Range("N95").Value = -35700
$ws.Range("H116").Value = 2066.6858
$ws.Range("I116").Value = 1914.96
$ws.Range("J116").Value = 2446
$ws.Range("K116").Value = 1914.96
$ws.Range("L116").Value = 2446
$ws.Range("M116").Value = 379.04
$ws.Range("N116").Value = -7034
$ws.Range("H136").Value = 1467
$ws.Range("I136").Value = 1188.2727
$ws.Range("K136").Value = 3564.8181
$ws.Range("M136").Value = -1014.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2066.6858
$ws.Range("I3").Value = 1914.96
$ws.Range("J3").Value = 2446
$ws.Range("K3").Value = 1914.96
$ws.Range("L3").Value = 2446
$ws.Range("M3").Value = -1800.96
$ws.Range("N3").Value = -2674
$ws.Range("H20").Value = 2716.8936
$ws.Range("J20").Value = 4576.6665
$ws.Range("L20").Value = 4576.6665
$ws.Range("N20").Value = -5070.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 555
$ws.Range("J11").Value = 555
$ws.Range("L11").Value = 555
$ws.Range("N11").Value = -835
$ws.Range("H31").Value = 1880.3182
$ws.Range("I31").Value = 1880.3182
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1880.3182
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1585.3182
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 1880.3182
$ws.Range("I34").Value = 1880.3182
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1880.3182
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1678.3182
$ws.Range("N34").ClearContents()
$ws.Range("H86").Value = 9330.333000000001
$ws.Range("I86").Value = 8999
$ws.Range("J86").Value = 9496
$ws.Range("K86").Value = 8999
$ws.Range("L86").Value = 9496
$ws.Range("M86").Value = -7876
$ws.Range("N86").Value = -11742
$ws.Range("H89").Value = 9330.333000000001
$ws.Range("I89").Value = 8999
$ws.Range("J89").Value = 9496
$ws.Range("K89").Value = 44995
$ws.Range("L89").Value = 47480
$ws.Range("M89").Value = -39379
$ws.Range("N89").Value = -58712
$ws.Range("H123").Value = 94807.57000000001
$ws.Range("J123").Value = 94807.57000000001
$ws.Range("L123").Value = 94807.57000000001
$ws.Range("N123").Value = -104607.57
$ws.Range("H132").Value = 3024.0667
$ws.Range("I132").Value = 2695.5833
$ws.Range("J132").Value = 4338
$ws.Range("K132").Value = 8086.749899999999
$ws.Range("L132").Value = 13014
$ws.Range("M132").Value = -5556.749899999999
$ws.Range("N132").Value = -18074

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 217.1
$ws.Range("I6").Value = 217.1
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 651.3
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -538.3
$ws.Range("N6").ClearContents()
$ws.Range("H11").Value = 100599.2
$ws.Range("I11").Value = 570.7143
$ws.Range("J11").Value = 333999
$ws.Range("K11").Value = 1712.1429
$ws.Range("L11").Value = 1001997
$ws.Range("M11").Value = -1572.1429
$ws.Range("N11").Value = -1002277
$ws.Range("H26").Value = 300.77777
$ws.Range("I26").Value = 184.66667
$ws.Range("J26").Value = 533
$ws.Range("K26").Value = 554.00001
$ws.Range("L26").Value = 1599
$ws.Range("M26").Value = -266.00001
$ws.Range("N26").Value = -2175
$ws.Range("H40").Value = 2908.5454
$ws.Range("I40").Value = 249.8
$ws.Range("J40").Value = 5124.1665
$ws.Range("K40").Value = 999.2
$ws.Range("L40").Value = 20496.666
$ws.Range("M40").Value = -930.2
$ws.Range("N40").Value = -20634.666
$ws.Range("H128").Value = 339556.38
$ws.Range("I128").Value = 339556.38
$ws.Range("K128").Value = 1018669.14
$ws.Range("M128").Value = -1013689.14
$ws.Range("H131").Value = 64834.938
$ws.Range("J131").Value = 2621.8333
$ws.Range("L131").Value = 7865.499899999999
$ws.Range("N131").Value = -17945.4999
$ws.Range("H138").Value = 4577.222
$ws.Range("I138").Value = 1476.6666
$ws.Range("J138").Value = 10778.333
$ws.Range("K138").Value = 4429.9998
$ws.Range("L138").Value = 32334.999
$ws.Range("M138").Value = 710.0002000000004
$ws.Range("N138").Value = -42614.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 36680
$ws.Range("J39").Value = 36680
$ws.Range("L39").Value = 36680
$ws.Range("N39").Value = -37744
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H122").Value = 2837.25
$ws.Range("I122").Value = 1675
$ws.Range("J122").Value = 3999.5
$ws.Range("K122").Value = 5025
$ws.Range("L122").Value = 11998.5
$ws.Range("M122").Value = -2575
$ws.Range("N122").Value = -16898.5
$ws.Range("H123").Value = 64479.5
$ws.Range("J123").Value = 64479.5
$ws.Range("L123").Value = 64479.5
$ws.Range("N123").Value = -69379.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5085.8335
$ws.Range("I40").Value = 3631.5
$ws.Range("J40").Value = 7994.5
$ws.Range("K40").Value = 3631.5
$ws.Range("L40").Value = 7994.5
$ws.Range("M40").Value = -3495.5
$ws.Range("N40").Value = -8266.5
$ws.Range("H61").Value = 2178.2
$ws.Range("I61").Value = 2178.2
$ws.Range("K61").Value = 2178.2
$ws.Range("M61").Value = -1976.2
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H113").Value = 2178.2
$ws.Range("I113").Value = 2178.2
$ws.Range("K113").Value = 2178.2
$ws.Range("M113").Value = -8.199999999999818
$ws.Range("H129").Value = 86000
$ws.Range("J129").Value = 86000
$ws.Range("L129").Value = 86000
$ws.Range("N129").Value = -96000
$ws.Range("H132").Value = 3326.5454
$ws.Range("I132").Value = 2262.75
$ws.Range("J132").Value = 6163.3335
$ws.Range("K132").Value = 6788.25
$ws.Range("L132").Value = 18490.0005
$ws.Range("M132").Value = -4258.25
$ws.Range("N132").Value = -23550.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H111").Value = 11248
$ws.Range("J111").Value = 11248
$ws.Range("L111").Value = 11248
